# Updated Quan's time log slide - ZWT
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# --- Paragraph 2: "Time Log (Starting 10/7/13)" -------------------------
# Split the trailing ")" into its own run (author re-typed the closing
# paren as a separate edit). Re-stamping the font size on the "prefix"
# sub-range forces PowerPoint to keep it as a distinct run even though
# the value does not change.
$para2 = $tr.Paragraphs(2, 1)
$prefix2 = $tr.Characters($para2.Start, $para2.Length - 2)
$prefix2.Font.Size = 25

# --- Paragraph 3: "10/07/13 – Discuss CImage versus Bitmap, ..." --------
# Insert two new "Prior - ..." bullet paragraphs ahead of the existing
# bullet, in reading order.
$para3 = $tr.Paragraphs(3, 1)
$para3.InsertBefore("Prior " + [char]0x2013 + " Research Sudoku fundamentals and algorithms for solving. Look into number creation." + [char]13)

$para4 = $tr.Paragraphs(4, 1)
$para4.InsertBefore("Prior " + [char]0x2013 + " Looked into GUI and MFC in C++" + [char]13)

# The original bullet is now paragraph 5. Split its leading "10/07/13 – "
# text into two runs: "10/07/13 " and "– ".
$para5 = $tr.Paragraphs(5, 1)
$prefix5 = $tr.Characters($para5.Start, 9)
$prefix5.Font.Size = 17
